$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''300.03'
$ws.Range("E2").Value = '''2.00%'
$ws.Range("D3").Value = '''32.27'
$ws.Range("E3").Value = '''3.65%'
$ws.Range("D4").Value = '''5.011'
$ws.Range("E4").Value = '''1.65%'
$ws.Range("D5").Value = '''0.07719'
$ws.Range("E5").Value = '''5.14%'
$ws.Range("D6").Value = '''2.267'
$ws.Range("E6").Value = '''-0.87%'
$ws.Range("D7").Value = '''7.938'
$ws.Range("E7").Value = '''2.53%'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '''0.9230'
$ws.Range("E8").Value = '''1.59%'
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = '''0.09911'
$ws.Range("E9").Value = '''23.80%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '''0.1762'
$ws.Range("E10").Value = '''4.29%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '''0.08449'
$ws.Range("E11").Value = '''4.61%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '''0.03307'
$ws.Range("E12").Value = '''6.77%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '''0.09872'
$ws.Range("E13").Value = '''-2.08%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '''0.001475'
$ws.Range("E14").Value = '''-2.76%'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '''0.005670'
$ws.Range("E15").Value = '''-1.51%'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '''3.544'
$ws.Range("E16").Value = '''1.74%'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '''3.822'
$ws.Range("E17").Value = '''1.95%'
$ws.Range("D18").Value = '''2.185'
$ws.Range("E18").Value = '''5.26%'
$ws.Range("D19").Value = '''0.3355'
$ws.Range("E19").Value = '''0.80%'
$ws.Range("D20").Value = '''0.1337'
$ws.Range("E20").Value = '''2.52%'
$ws.Range("D21").Value = '''4.125'
$ws.Range("E21").Value = '''3.64%'
$ws.Range("D22").Value = '''0.2085'
$ws.Range("E22").Value = '''-0.79%'
$ws.Range("D23").Value = '''0.04534'
$ws.Range("E23").Value = '''-0.29%'
$ws.Range("D24").Value = '''0.001217'
$ws.Range("E24").Value = '''0.55%'
$ws.Range("D25").Value = '''0.004369'
$ws.Range("E25").Value = '''-5.99%'
$ws.Range("D26").Value = '''0.0001292'
$ws.Range("E26").Value = '''-0.73%'
$ws.Range("D27").Value = '''0.0003375'
$ws.Range("E27").Value = '''-0.66%'
$ws.Range("D39").Value = '''0.01702'
$ws.Range("E39").Value = '''6.09%'
$ws.Range("D40").Value = '''0.04676'
$ws.Range("E40").Value = '''5.01%'
$ws.Range("D41").Value = '''0.007717'
$ws.Range("E41").Value = '''5.43%'
$ws.Range("D42").Value = '''0.009775'
$ws.Range("E42").Value = '''13.06%'
$ws.Range("D43").Value = '''0.1396'
$ws.Range("E43").Value = '''4.93%'
$ws.Range("D44").Value = '''0.002101'
$ws.Range("E44").Value = '''7.88%'
$ws.Range("D45").Value = '''0.009702'
$ws.Range("E45").Value = '''1.84%'
$ws.Range("D46").Value = '''0.00006066'
$ws.Range("E46").Value = '''1.92%'
$ws.Range("D47").Value = '''0.00000000746'
$ws.Range("E47").Value = '''-0.66%'
$ws.Range("D48").Value = '''2.551'
$ws.Range("E48").Value = '''13.83%'
$ws.Range("D49").Value = '''0.001988'
$ws.Range("E49").Value = '''-31.40%'
$ws.Range("D50").Value = '''0.00002088'
$ws.Range("E50").Value = '''-0.66%'
$ws.Range("D51").Value = '''0.0001988'
$ws.Range("E51").Value = '''-0.66%'

Write-Output "Applied cryptos update"
